# Apply the "Geomorphic Potential" -> "Geomorphic Potential (Confinement)" update
# and add Habitat_Type "Confinement" + Data_Sources "USFS Valley Confinement Algorithm "
# to rows 26-28, plus highlight the new Habitat_Type cells and turn on AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column A (Habitat_Quality_Scoring_Metric) label for rows 26-28
$ws.Range("A26:A28").Value = "Geomorphic Potential (Confinement)"

# Populate column B (Habitat_Type) for rows 26-28 with "Confinement"
$ws.Range("B26:B28").Value = "Confinement"

# Populate column D (Data_Sources) for rows 26-28
$ws.Range("D26:D28").Value = "USFS Valley Confinement Algorithm "

# Highlight the new Habitat_Type cells with a light fill (theme color 9, light tint)
$ws.Range("B26:B28").Interior.ThemeColor = 9
$ws.Range("B26:B28").Interior.TintAndShade = 0.8

# Turn on AutoFilter for the data range
$ws.Range("A1:L1").AutoFilter() | Out-Null

# Adjust the view similarly to the authored workbook
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("D28").Select()
